$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 6485
$ws.Range("J58").Value = 16000
$ws.Range("L58").Value = 48000
$ws.Range("N58").Value = -48300
$ws.Range("H69").Value = 7015
$ws.Range("J69").Value = 7015
$ws.Range("L69").Value = 21045
$ws.Range("N69").Value = -22793
$ws.Range("H72").Value = 7015
$ws.Range("J72").Value = 7015
$ws.Range("L72").Value = 63135
$ws.Range("N72").Value = -71871
$ws.Range("H116").Value = 21178.285
$ws.Range("I116").Value = 5142.4287
$ws.Range("K116").Value = 5142.4287
$ws.Range("M116").Value = -1700.4287
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 619
$ws.Range("I4").Value = 388.5
$ws.Range("K4").Value = 388.5
$ws.Range("M4").Value = -272.5
$ws.Range("H30").Value = 4833.3335
$ws.Range("I30").Value = 4750
$ws.Range("K30").Value = 4750
$ws.Range("M30").Value = -4600
$ws.Range("H32").Value = 255904.7
$ws.Range("I32").Value = 372785.34
$ws.Range("J32").Value = 13152.615
$ws.Range("K32").Value = 372785.34
$ws.Range("L32").Value = 13152.615
$ws.Range("M32").Value = -372498.34
$ws.Range("N32").Value = -13726.615
$ws.Range("H132").Value = 545086.1
$ws.Range("I132").Value = 557174.7
$ws.Range("J132").Value = 1099
$ws.Range("K132").Value = 1671524.1
$ws.Range("L132").Value = 3297
$ws.Range("M132").Value = -1668994.1
$ws.Range("N132").Value = -8357
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 450
$ws.Range("J8").Value = 500
$ws.Range("L8").Value = 500
$ws.Range("N8").Value = -780
$ws.Range("H11").Value = 116.181816
$ws.Range("J11").Value = 280.66666
$ws.Range("L11").Value = 280.66666
$ws.Range("N11").Value = -560.66666
$ws.Range("H12").Value = 362.33334
$ws.Range("I12").Value = 280
$ws.Range("J12").Value = 527
$ws.Range("K12").Value = 280
$ws.Range("L12").Value = 527
$ws.Range("M12").Value = -112
$ws.Range("N12").Value = -863
$ws.Range("H86").Value = 2153.484
$ws.Range("I86").Value = 2096.0557
$ws.Range("K86").Value = 2096.0557
$ws.Range("M86").Value = -973.0556999999999
$ws.Range("H89").Value = 2153.484
$ws.Range("I89").Value = 2096.0557
$ws.Range("K89").Value = 10480.2785
$ws.Range("M89").Value = -4864.2785
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 4499.3335
$ws.Range("I21").Value = 4499
$ws.Range("J21").Value = 4500
$ws.Range("K21").Value = 4499
$ws.Range("L21").Value = 4500
$ws.Range("N21").Value = -4970
$ws.Range("M21").Value = -4264
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H47").Value = 24000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 24000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 24000
$ws.Range("N47").Value = -25132
$ws.Range("M47").ClearContents()
$ws.Range("H62").Value = 3378
$ws.Range("I62").Value = 3199.75
$ws.Range("J62").Value = 3734.5
$ws.Range("K62").Value = 3199.75
$ws.Range("L62").Value = 3734.5
$ws.Range("M62").Value = -2575.75
$ws.Range("N62").Value = -4982.5
$ws.Range("H65").Value = 3378
$ws.Range("I65").Value = 3199.75
$ws.Range("J65").Value = 3734.5
$ws.Range("K65").Value = 15998.75
$ws.Range("L65").Value = 18672.5
$ws.Range("M65").Value = -12878.75
$ws.Range("N65").Value = -24912.5
$ws.Range("H95").Value = 33194
$ws.Range("J95").Value = 33194
$ws.Range("L95").Value = 33194
$ws.Range("N95").Value = -38686
$ws.Range("H96").Value = 7421.4287
$ws.Range("J96").Value = 7421.4287
$ws.Range("L96").Value = 7421.4287
$ws.Range("N96").Value = -12913.4287
$ws.Range("H107").Value = 515.6923
$ws.Range("I107").Value = 495
$ws.Range("K107").Value = 495
$ws.Range("M107").Value = 1425
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4996.5
$ws.Range("I21").Value = 4994.75
$ws.Range("K21").Value = 4994.75
$ws.Range("M21").Value = -4821.75
$ws.Range("H30").Value = 4996.5
$ws.Range("I30").Value = 4994.75
$ws.Range("K30").Value = 4994.75
$ws.Range("M30").Value = -4889.75
$ws.Range("H80").Value = 4384.36
$ws.Range("I80").Value = 3176.4119
$ws.Range("J80").Value = 6951.25
$ws.Range("K80").Value = 3176.4119
$ws.Range("L80").Value = 6951.25
$ws.Range("M80").Value = -2178.4119
$ws.Range("N80").Value = -8947.25
$ws.Range("H83").Value = 4384.36
$ws.Range("I83").Value = 3176.4119
$ws.Range("J83").Value = 6951.25
$ws.Range("K83").Value = 15882.0595
$ws.Range("L83").Value = 34756.25
$ws.Range("M83").Value = -10890.0595
$ws.Range("N83").Value = -44740.25
$ws.Range("H122").Value = 4289.625
$ws.Range("I122").Value = 3515.7058
$ws.Range("K122").Value = 10547.1174
$ws.Range("M122").Value = -8097.117400000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 5845
$ws.Range("I9").Value = 293.33334
$ws.Range("J9").Value = 22500
$ws.Range("K9").Value = 293.33334
$ws.Range("L9").Value = 22500
$ws.Range("M9").Value = -69.33334000000002
$ws.Range("N9").Value = -22948
$ws.Range("H18").Value = 15000
$ws.Range("J18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("N18").Value = -15344
$ws.Range("H34").Value = 14008
$ws.Range("J34").Value = 14008
$ws.Range("L34").Value = 14008
$ws.Range("N34").Value = -14352
$ws.Range("H132").Value = 5083140
$ws.Range("I132").Value = 10624645
$ws.Range("J132").Value = 3427.0833
$ws.Range("K132").Value = 31873935
$ws.Range("L132").Value = 10281.2499
$ws.Range("M132").Value = -31871405
$ws.Range("N132").Value = -15341.2499
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 60000
$ws.Range("I70").Value = 60000
$ws.Range("K70").Value = 60000
$ws.Range("M70").Value = -59685
$ws.Range("H73").Value = 60000
$ws.Range("I73").Value = 60000
$ws.Range("K73").Value = 60000
$ws.Range("M73").Value = -58908
$ws.Range("H107").Value = 5493.421
$ws.Range("I107").Value = 5020.0713
$ws.Range("J107").Value = 6818.8
$ws.Range("K107").Value = 15060.2139
$ws.Range("L107").Value = 20456.4
$ws.Range("M107").Value = -13140.2139
$ws.Range("N107").Value = -24296.4
$ws.Range("H126").Value = 3038.8333
$ws.Range("I126").Value = 3014.3572
$ws.Range("J126").Value = 3124.5
$ws.Range("K126").Value = 9043.0716
$ws.Range("L126").Value = 9373.5
$ws.Range("M126").Value = -6573.071599999999
$ws.Range("N126").Value = -14313.5
$ws.Range("H132").Value = 5557700.5
$ws.Range("I132").Value = 5954372
$ws.Range("K132").Value = 17863116
$ws.Range("M132").Value = -17860586
$ws.Range("H136").Value = 24348906
$ws.Range("I136").Value = 5436020
$ws.Range("J136").Value = 100000450
$ws.Range("K136").Value = 16308060
$ws.Range("L136").Value = 300001350
$ws.Range("M136").Value = -16305510
$ws.Range("N136").Value = -300006450
